# Adds two new LeetCode entries to the tracker sheet:
#   39. Combination Sum  (row 11)
#   78. Subsets          (row 12)
# Both are tagged "Java" and get a "date solved" value in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$xlTop    = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignTop

# ---- Row 11: 39. Combination Sum ----
$ws.Range("A11").Value = 39
$ws.Range("A11").HorizontalAlignment = $xlCenter
$ws.Range("A11").VerticalAlignment = $xlTop

$ws.Range("B11").Value = "Combination Sum"
$ws.Range("B11").Style = "Normal"

$ws.Range("C11").Value = "Java"

$ws.Range("D11").Value = 44982
$ws.Range("D11").NumberFormat = "d-mmm-yy"

# ---- Row 12: 78. Subsets ----
$ws.Range("A12").Value = 78
$ws.Range("A12").HorizontalAlignment = $xlCenter
$ws.Range("A12").VerticalAlignment = $xlTop

$ws.Range("B12").Value = "Subsets"
$ws.Range("B12").Style = "Normal"

$ws.Range("C12").Value = "Java"

$ws.Range("D12").Value = 44981
$ws.Range("D12").NumberFormat = "d-mmm-yy"

# Autosize the new date column, like Excel would after typing dates in it.
$ws.Columns("D").AutoFit()

# Matches the cursor position left by the author in the saved workbook.
$ws.Range("B17").Select()
